# Updates Betfair back/lay odds data (columns F:AO) for rows 2-12 on Sheet1,
# matching the values captured in the "Jogos_do_Dia_Betfair_Back_Lay_2025-10-15" export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colCount = 36

$values = @(1.8, 1.99, 4, 4.9, 3.75, 4.5, 1.33, 1.05, 4.7, 1.22, 2.24, 1.68, 1.53, 2.64, 1.63, 2.34, 1.26, 2.02, 23, 1000, 100, 1000, 14.5, 10, 18, 48, 13.5, 10.5, 18, 55, 23, 18, 29, 80, 10, 42)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F2:AO2").Value = $rowData

$values = @(1.41, 1.42, 9.4, 9.6, 5.3, 5.5, 1.33, 1.04, 5, 1.24, 2.34, 1.72, 1.53, 2.84, 2.04, 1.92, 1.11, 3.4, 22, 34, 85, 320, 8.800000000000001, 12, 32, 140, 8.199999999999999, 9.800000000000001, 26, 120, 11.5, 14.5, 36, 140, 6, 170)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F3:AO3").Value = $rowData

$values = @(4.8, 4.9, 2, 2.02, 3.35, 3.4, 1.54, 1.12, 3, 1.48, 1.66, 2.46, 1.24, 4.9, 2.16, 1.84, 1.98, 1.25, 9.199999999999999, 7, 11, 23, 13, 7.8, 10.5, 25, 34, 19.5, 23, 50, 120, 75, 95, 170, 110, 21)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F4:AO4").Value = $rowData

$values = @(2.18, 2.2, 4, 4.2, 3.35, 3.4, 1.49, 1.09, 3.3, 1.41, 1.78, 2.24, 1.29, 4.2, 1.94, 2.02, 1.32, 1.83, 11.5, 13, 28, 85, 8.199999999999999, 7.4, 16.5, 55, 13, 10.5, 19, 70, 27, 24, 44, 130, 21, 75)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F5:AO5").Value = $rowData

$values = @(2.56, 2.6, 3.25, 3.35, 3.25, 3.3, 1.55, 1.11, 2.9, 1.51, 1.62, 2.54, 1.23, 5.2, 2.08, 1.89, 1.43, 1.63, 9.4, 9.800000000000001, 19.5, 60, 8.199999999999999, 7, 14.5, 48, 15, 12, 22, 70, 38, 34, 60, 160, 42, 60)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F6:AO6").Value = $rowData

$values = @(2.98, 3.1, 2.74, 2.82, 3.15, 3.25, 1.53, 1.11, 3, 1.47, 1.67, 2.44, 1.24, 4.8, 1.97, 1.92, 1.54, 1.48, 10, 9.199999999999999, 16.5, 44, 9.6, 7.2, 13, 36, 18.5, 13, 21, 1000, 55, 40, 1000, 160, 46, 42)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F7:AO7").Value = $rowData

$values = @(2.18, 2.2, 4.1, 4.2, 3.25, 3.3, 1.57, 1.13, 2.82, 1.53, 1.6, 2.6, 1.21, 5.3, 2.22, 1.78, 1.31, 1.83, 8.6, 11.5, 27, 95, 7.2, 7.6, 18, 70, 11.5, 11.5, 24, 90, 26, 30, 60, 180, 29, 95)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F8:AO8").Value = $rowData

$values = @(2.44, 2.5, 3.25, 3.35, 3.4, 3.5, 1.46, 1.09, 3.5, 1.38, 1.84, 2.14, 1.31, 4, 1.89, 2.06, 1.42, 1.66, 12, 11.5, 21, 60, 9.4, 7.4, 14, 40, 15, 11.5, 18.5, 55, 34, 28, 44, 110, 24, 42)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F9:AO9").Value = $rowData

$values = @(3.6, 3.7, 2.4, 2.44, 3.15, 3.2, 1.58, 1.13, 2.72, 1.56, 1.57, 2.7, 1.19, 5.7, 2.22, 1.8, 1.69, 1.37, 8.199999999999999, 7.4, 13, 36, 10, 7.4, 12, 34, 23, 16, 24, 65, 75, 60, 85, 180, 90, 34)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F10:AO10").Value = $rowData

$values = @(2.24, 2.6, 3.2, 4.1, 3, 3.55, 1.43, 1.08, 3.4, 1.34, 1.84, 2, 1.31, 3.55, 1.75, 2, 1.34, 1.64, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000, 1000)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F11:AO11").Value = $rowData

$values = @(2.24, 2.26, 3.3, 3.45, 3.8, 4, 1.44, 1.07, 3.65, 1.34, 1.93, 2.02, 1.35, 3.6, 1.79, 2.06, 1.4, 1.79, 15, 13.5, 25, 65, 9.800000000000001, 8.199999999999999, 14, 55, 14.5, 11, 19.5, 55, 30, 24, 40, 100, 18.5, 44)
$rowData = New-Object "object[,]" 1, $colCount
for ($i = 0; $i -lt $colCount; $i++) { $rowData[0, $i] = [double]$values[$i] }
$ws.Range("F12:AO12").Value = $rowData

Write-Output "Updated odds values in F2:AO12."
